$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Starting deck (4 slides):
#   1: Textbook Rental Library / subtitle          (unchanged)
#   2: Use Cases / UC1, UC2
#   3: TRLApp / Show running code
#   4: Code Smells / ...                            (unchanged)
#
# Target deck (6 slides):
#   1: Textbook Rental Library / subtitle          (unchanged)
#   2: DoCD / (empty)                               NEW slide
#   3: Use Cases / UC1, UC2                         (same content as old slide 2)
#   4: Demonstration: Core / Out, In, JUnit Tests & Coverage, Software Class Diagram, ...
#   5: Demonstration: Additional Features / ...      NEW slide
#   6: Code Smells / ...                             (unchanged)
# ---------------------------------------------------------------------------

# --- Step 1: duplicate slide 2 ("Use Cases") -------------------------------
# The duplicate lands immediately after slide 2 (i.e. at position 3) and
# keeps the original "Use Cases" / "UC1 / UC2" content untouched. We then
# turn the original slide 2 into the new "DoCD" slide.
$s2 = $p.Slides.Item(2)
$dup = $s2.Duplicate()

# --- Step 2: turn (original) slide 2 into the new "DoCD" slide -------------
$titleTr = $s2.Shapes.Item(1).TextFrame.TextRange
$titleTr.Delete()
$titleTr.Text = "DoCD"

$bodyTr = $s2.Shapes.Item(2).TextFrame.TextRange
$bodyTr.Delete()

# --- Step 3: slide 3 ("Use Cases") is the untouched duplicate, no edits ----

# --- Step 4: edit slide 4 ("TRLApp") into "Demonstration: Core" -----------
$s4 = $p.Slides.Item(4)

$titleTr4 = $s4.Shapes.Item(1).TextFrame.TextRange
$titleTr4.Delete()
$titleTr4.Text = "Demonstration: Core"

$bodyTr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$bodyTr4.Delete()
$bodyTr4.Text = "Out`rIn`rJUnit Tests & Coverage`rSoftware Class Diagram`r`r`r"
$lastPara = $bodyTr4.Paragraphs(7, 1)
$lastPara.IndentLevel = 2
$lastPara.Font.Size = 22

# --- Step 5: duplicate slide 4 to create slide 5 ---------------------------
$dup2 = $s4.Duplicate()
$s5 = $p.Slides.Item(5)

$titleTr5 = $s5.Shapes.Item(1).TextFrame.TextRange
$titleTr5.Delete()
$titleTr5.Text = "Demonstration: Additional Features"

$bodyTr5 = $s5.Shapes.Item(2).TextFrame.TextRange
$bodyTr5.Delete()
$bodyTr5.Text = "Event logging of changes to Patron and Copy state`rSearching through past Event logs`rCreate and attach Holds to all Patrons with overdue Copies`rGenerate overdue notices for all Patrons with overdue Holds`rPrint overdue notices for all Patrons with overdue Holds"

# --- Step 6: slide 6 ("Code Smells") is the original slide 4, untouched ---

Write-Host "Final slide count:" $p.Slides.Count
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    Write-Host "Slide $i title:" $s.Shapes.Item(1).TextFrame.TextRange.Text
}
